$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-10-08 to 2023-10-09
# (Excel serial date 45207 -> 45208)
$newDate = Get-Date -Year 2023 -Month 10 -Day 9 -Hour 0 -Minute 0 -Second 0
$newDate = $newDate.Date
$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
